# Automatische test-sync: 2025-06-18 16:30:10
# Append the new "Afmelding nieuwsbrief" log row (row 18) to the Logs sheet,
# extend the conditional-formatting ranges to include it, and bump the
# "Afmelding" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row of data
$logs.Range("A18").Value = "Afmelding nieuwsbrief"
$logs.Range("B18").Value = "mailmind.test@zohomail.eu"
$logs.Range("C18").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D18").Value = "Afmelding"
$logs.Range("F18").Value = "2025-06-18 16:00:11"
$logs.Range("G18").Value = "Nee"

# Extend the conditional formatting ranges so the new row is covered too
$dFormatConditions = $logs.Range("D2:D17").FormatConditions
for ($i = 1; $i -le $dFormatConditions.Count(); $i++) {
    $dFormatConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D18"))
}

$gFormatConditions = $logs.Range("G2:G17").FormatConditions
for ($i = 1; $i -le $gFormatConditions.Count(); $i++) {
    $gFormatConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G18"))
}

# Update the Dashboard summary count for "Afmelding" (2 -> 3)
$dashboard.Range("B4").Value = 3
